# "Generate Report for Archive"
#
# 1. Every cell whose status text was "Ready for handoff" becomes "In Translation"
#    (Overview!E2:F2, E3:F3 -- one column per locale -- and the per-locale detail
#    sheets' "Status" column, zh-cn!C2:C3 / de-de!C2:C3).
# 2. The two "Status"-ish columns that were sized for the old, longer text
#    ("Ready for handoff") are narrowed to fit the shorter "In Translation" text:
#    Overview columns E and F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# Narrower columns to match the shorter "In Translation" status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
